$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 162, shifting existing rows 162:183 down to 163:184.
$ws.Rows.Item(162).Insert()

# Copy the row that now sits at 163 (the former row 162, preserved by the
# insert/shift) into the new blank row 162 so that formatting and all the
# values that stay constant across the two records are already correct.
$ws.Rows.Item(163).Copy()
$ws.Rows.Item(162).PasteSpecial()

# Now overwrite just the handful of fields that differ for the new record.
$ws.Cells.Item(162, 4).Value = 45142
$ws.Cells.Item(162, 13).Value = 320
$ws.Cells.Item(162, 14).Value = 17000
$ws.Cells.Item(162, 15).Value = 18000
$ws.Cells.Item(162, 16).Value = 17469
$ws.Cells.Item(162, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(162, 19).Value = 873
